$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped from
# 45182 (2023-09-13) to 45184 (2023-09-15) for every data row (rows 2-270).
$newDate = (Get-Date -Year 2023 -Month 9 -Day 15).Date

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 270 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = $newDate
